# Generate Report for Archive
#
# 1) The shared status text "Ready for handoff" becomes "In Translation"
#    everywhere it is used: the two language-status columns (zh-cn, de-de)
#    on the "Overview" sheet, and the "Status" column on each of the
#    per-language report sheets ("zh-cn", "de-de").
# 2) The corresponding status columns are narrowed (report columns no
#    longer need to be as wide once the longer "Ready for handoff" text
#    is replaced with the shorter "In Translation").

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newWidth = 12.5   # -> stored column width ~13.33 (closest grid value to the narrower target)

# --- Overview sheet: zh-cn status (col E) and de-de status (col F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet: Status column (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
